$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36, shifting existing rows 36-74 down to 37-75
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new record
$ws.Range("A36").Value = 6
$ws.Range("B36").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C36").Value = "Metropolitana"
$ws.Range("D36").Value = 44669
$ws.Range("D36").NumberFormat = $ws.Range("D37").NumberFormat
$ws.Range("E36").Value = 13
$ws.Range("F36").Value = 100114007
$ws.Range("G36").Value = "Jengibre"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 230
$ws.Range("K36").Value = 10000
$ws.Range("L36").Value = 11000
$ws.Range("M36").Value = 10348
$ws.Range("N36").Value = "`$/caja 13 kilos"
$ws.Range("O36").Value = "Perú"
$ws.Range("P36").Value = 796
$ws.Range("Q36").Value = 13
$ws.Range("R36").Value = "Hortaliza"
